$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FamilyData")

# Update existing photo file names to the new "id<PersonID>-name.png" convention
$ws.Range("Q2").Value  = "id1-leonid.png"
$ws.Range("Q3").Value  = "id2-elena.png"
$ws.Range("Q4").Value  = "id3-sergey.png"
$ws.Range("Q6").Value  = "id5-maksim.png"
$ws.Range("Q7").Value  = "id6-sophia.png"
$ws.Range("Q12").Value = "id11-milana.png"

# Fill in previously empty photo cells
$ws.Range("Q8").Value  = "id7-sviatik.png"
$ws.Range("Q30").Value = "id29-lilina.png"
$ws.Range("Q32").Value = "id31-sofia.png"
$ws.Range("Q33").Value = "id32-vasilii.png"
$ws.Range("Q34").Value = "id33-fedor.png"
$ws.Range("Q35").Value = "id34-lelia.png"
$ws.Range("Q36").Value = "id35-nikolay.png"
$ws.Range("Q39").Value = "id38-ivan.png"
$ws.Range("Q52").Value = "id51-kazimir.png"
$ws.Range("Q53").Value = "id52-fedor.png"

# Fill in a previously empty Occupation note
$ws.Range("P52").Value = "Фотография Каземира сделана  в г. Ровны в 1928 году"
